$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# --- Sheet1: header row becomes "Figure" + re-pointed labels, Figure column becomes numeric ---
$ws1.Range("A1").Value = "Figure"
$ws1.Range("B1").Value = "Square Size"
$ws1.Range("C1").Value = "Origin X"
$ws1.Range("D1").Value = "Origin Y"
$ws1.Range("E1").Value = "Image H"
$ws1.Range("F1").Value = "Image W"

$ws1.Range("A2").Value = 1
$ws1.Range("A3").Value = 2
$ws1.Range("A4").Value = 4
$ws1.Range("A5").Value = 5
$ws1.Range("A6").Value = 6
$ws1.Range("A7").Value = 7
$ws1.Range("A8").Value = 8
$ws1.Range("A9").Value = 9
$ws1.Range("A10").Value = 10

# --- Sheet2: re-write the "+/-" label cell so the shared-string table is rebuilt ---
# (leading apostrophe keeps the cell's existing quote-prefixed text style, s="1")
$ws2.Range("B6").Value = "'+/-"

# --- Sheet1: page setup now records an explicit portrait orientation ---
$ws1.PageSetup.Orientation = 1

# --- View state: Sheet1 becomes the selected/active tab with F12:G16 selected ---
$ws1.Activate()
$ws1.Range("F12:G16").Select()
